$d = $word.ActiveDocument

# CAI number 1 -> 4
$d.Content.Find.Execute("CAI N°1", $true, $false, $false, $false, $false, $true, 1, $false, "CAI N°4", 2) | Out-Null
# C.cerrado -> C.abierto
$d.Content.Find.Execute("C.cerrado", $true, $false, $false, $false, $false, $true, 1, $false, "C.abierto", 2) | Out-Null
# invoice date (must run before Fecha de Inicio change)
$d.Content.Find.Execute("2022-01-09", $true, $false, $false, $false, $false, $true, 1, $false, "2022-01-10", 2) | Out-Null
# client name
$d.Content.Find.Execute("PC INC.", $true, $false, $false, $false, $false, $true, 1, $false, "CENCOSUD S.A.", 2) | Out-Null
# client RUT
$d.Content.Find.Execute("7794143-6", $true, $false, $false, $false, $false, $true, 1, $false, "12621140-6", 2) | Out-Null
# address / comuna
$d.Content.Find.Execute("/la florida", $true, $false, $false, $false, $false, $true, 1, $false, "/maipu", 2) | Out-Null
# vencimiento date
$d.Content.Find.Execute("21/12/2021", $true, $false, $false, $false, $false, $true, 1, $false, "09/02/2022", 2) | Out-Null
# course name
$d.Content.Find.Execute("Aplicación De Herramientas De Redacción Y Ortografía Para Profesionales", $true, $false, $false, $false, $false, $true, 1, $false, "Diseño De Proyectos A Través De La Metodología De Design Thinking", 2) | Out-Null
# codigo sence
$d.Content.Find.Execute("1238020242", $true, $false, $false, $false, $false, $true, 1, $false, "1237971736", 2) | Out-Null
# Nº Horas value
$d.Content.Find.Execute("Nº Horas: 42", $true, $false, $false, $false, $false, $true, 1, $false, "Nº Horas: 16", 2) | Out-Null
# Fecha de Inicio (must run after invoice date change)
$d.Content.Find.Execute(" 2021-11-20", $true, $false, $false, $false, $false, $true, 1, $false, " 2022-01-09", 2) | Out-Null
# Fecha de Termino
$d.Content.Find.Execute("2021-10-19", $true, $false, $false, $false, $false, $true, 1, $false, "2022-06-09", 2) | Out-Null
# Nº Registro Sence value
$d.Content.Find.Execute("Nº Registro Sence: 1", $true, $false, $false, $false, $false, $true, 1, $false, "Nº Registro Sence: 133", 2) | Out-Null
# amount, 3 occurrences all change
$d.Content.Find.Execute("168000", $true, $false, $false, $false, $false, $true, 1, $false, "80000", 2) | Out-Null
# empty OTRO value -> None
$d.Content.Find.Execute("OTRO: ", $true, $false, $false, $false, $false, $true, 1, $false, "OTRO: None", 2) | Out-Null
# Orden de Compra N°
$d.Content.Find.Execute("12313", $true, $false, $false, $false, $false, $true, 1, $false, "12345", 2) | Out-Null
# Obs text
$d.Content.Find.Execute("ivan weco", $true, $false, $false, $false, $false, $true, 1, $false, "Prueba query nueva", 2) | Out-Null

Write-Output "done"
